$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers refreshed)
$wb.Worksheets.Item(1).Name = "GNG_TO-1650996188884854"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961907008893"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961907008893"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961907728522"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961908368897"

# Sheet 1: GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961888528545.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961888688567.csv"
$ws1.Range("B4").Value = "go_stims-16509961888688567.csv"
$ws1.Range("B5").Value = "GNG_stims-1650996188884854.csv"

# Sheet 2: NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650996190172855.csv"
$ws2.Range("B3").Value = "TB-1650996190268862.csv"
$ws2.Range("B4").Value = "OB-1650996189724863.csv"
$ws2.Range("B5").Value = "ZB-match_0-16509961896688988.csv"
$ws2.Range("B6").Value = "TB-16509961905328608.csv"
$ws2.Range("B7").Value = "ZB-match_0-16509961895408885.csv"
$ws2.Range("B8").Value = "TB-1650996190684854.csv"
$ws2.Range("B9").Value = "OB-16509961901328912.csv"
$ws2.Range("B10").Value = "ZB-match_5-16509961893888528.csv"

# Sheet 4: TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961907248576.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961907088575.csv"
$ws4.Range("B4").Value = "MM_stims-16509961907488916.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961907248576.csv"
$ws4.Range("B6").Value = "MM_stims-16509961907728522.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961907488916.csv"

# Sheet 5: vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16509961907728522.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961908208907.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961908048604.csv"
$ws5.Range("B5").Value = "SAT_stims-16509961907888951.csv"
